# The document contains five "<id>...</id>" tags (p104r_1 .. p104r_5), each
# split across three separate runs: "<id>", the bare id text (plain black
# run), and "</id>". Collapse each trio into a single run so the full tag
# text lives in one run (taking on the Courier-New / gold-colored "<id>"
# run's formatting), matching how the tag text is represented elsewhere in
# the document (e.g. "<head>...", "<div>...").

$d = $word.ActiveDocument

for ($n = 1; $n -le 5; $n++) {
    $tag = "<id>p104r_$n</id>"
    $found = $d.Content.Find.Execute($tag, $true, $false, $false, $false, $false, $true, 1, $false, $tag, 2)
    Write-Output "p104r_$n merged: $found"
}
